$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1748
$ws1.Range("F5").Value = 774
$ws1.Range("F6").Value = 84

# Sheet "全部类型" (all types / combined list)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1748
$ws4.Range("F6").Value = 774
$ws4.Range("F7").Value = 84
